$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the string "admin123" in B2 with the numeric value 123
$ws.Range("B2").Value = 123

# Update the active selection to B3
$ws.Range("B3").Select()
